$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (now "B")
$ws.Range("A2").Value = "B"
$ws.Range("B2").Value = 0.9256756756756757
$ws.Range("C2").Value = 0.958041958041958
$ws.Range("D2").Value = 0.9415807560137457
$ws.Range("E2").Value = 143

# Row 3 (now "M")
$ws.Range("A3").Value = "M"
$ws.Range("B3").Value = 0.925
$ws.Range("C3").Value = 0.8705882352941177
$ws.Range("D3").Value = 0.896969696969697
$ws.Range("E3").Value = 85

# Row 4 (accuracy)
$ws.Range("B4").Value = 0.9254385964912281
$ws.Range("C4").Value = 0.9254385964912281
$ws.Range("D4").Value = 0.9254385964912281
$ws.Range("E4").Value = 0.9254385964912281

# Row 5 (macro avg)
$ws.Range("B5").Value = 0.9253378378378379
$ws.Range("C5").Value = 0.9143150966680378
$ws.Range("D5").Value = 0.9192752264917213
$ws.Range("E5").Value = 228

# Row 6 (weighted avg)
$ws.Range("B6").Value = 0.9254237790422001
$ws.Range("C6").Value = 0.9254385964912281
$ws.Range("D6").Value = 0.9249494401420609
$ws.Range("E6").Value = 228
